# Daily attendance processing - swap the order of "Recorded By" entries
# from "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every row in the "Recorded By" column (G) that currently has the
# old ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$firstRow = $used.Row

$oldVal = "System, dnasr281@gmail.com"
$newVal = "dnasr281@gmail.com, System"

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstRow + $i
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldVal) {
        $cell.Value = $newVal
    }
}
